$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.091.67"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "3.177.34"
$ws.Range("E3").Value = "  -3.74%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.81"
$ws.Range("E5").Value = "  -2.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.90"
$ws.Range("E6").Value = "  -4.58%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.173.69"
$ws.Range("E8").Value = "  -3.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("E10").Value = "  -5.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("E11").Value = "  -6.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -3.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -4.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.31"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "3.697.75"
$ws.Range("E15").Value = "  -3.96%  "

$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "3.173.12"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("D18").Value = "63.045.21"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("E19").Value = "  -4.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.61"
$ws.Range("E20").Value = "  -4.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  -5.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("E23").Value = "  -5.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").Value = "  -4.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.36"
$ws.Range("E25").Value = "  -3.28%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  -4.42%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("E29").Value = "  -6.43%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.71"
$ws.Range("E30").Value = "  -7.43%  "

$ws.Range("E31").Value = "  -5.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.21"
$ws.Range("E32").Value = "  -5.09%  "

$ws.Range("E33").Value = "  -4.29%  "

$ws.Range("E34").Value = "  -6.74%  "

$ws.Range("E35").Value = "  -5.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("E36").Value = "  -4.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.09"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("D38").Value = "0.0₃0707"
$ws.Range("E38").Value = "  -4.85%  "

$ws.Range("E39").Value = "  -3.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "401.01"
$ws.Range("E40").Value = "  -7.04%  "

$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  -3.31%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.811.88"
$ws.Range("E43").Value = "  -9.74%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.111"
$ws.Range("E44").Value = "  -6.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  -5.08%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.02"
$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.53"
$ws.Range("E49").Value = "  -5.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.08"
$ws.Range("E50").Value = "  -4.76%  "

$ws.Range("E51").Value = "  -2.31%  "
